# Add dashboard logic and fix input type on sales view
#
# The "Fecha Venta" (sale date) column stores values as plain text
# (e.g. "2024-10-29") rather than real Excel dates. Assigning a
# date-shaped string straight to .Value lets Excel auto-convert it to a
# date serial, so we briefly force the cell to Text format, assign the
# string, then restore the style back to Normal (keeps the text value,
# drops the number format override).
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: ID Venta=1, XYZ-456, 2024-10-29, 7835, 3, 2
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "XYZ-456"
Set-TextValue $ws.Range("C2") "2024-10-29"
$ws.Range("D2").Value = 7835
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 2

# New row: ID Venta=2, HOW-349, 2024-10-29, 7835, 3, 3
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "HOW-349"
Set-TextValue $ws.Range("C3") "2024-10-29"
$ws.Range("D3").Value = 7835
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 3

# ID Venta=3, ABC003, 2022-04-01, 19999.75, 1, 1 (price corrected from 2500.64)
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ABC003"
Set-TextValue $ws.Range("C4") "2022-04-01"
$ws.Range("D4").Value = 19999.75
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1

# ID Venta=4, ABC011, 2023-04-10, 27999.99, 2, 2
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "ABC011"
Set-TextValue $ws.Range("C5") "2023-04-10"
$ws.Range("D5").Value = 27999.99
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 2

# ID Venta=5, ABC015, 2023-05-15, 23999, 3, 3
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "ABC015"
Set-TextValue $ws.Range("C6") "2023-05-15"
$ws.Range("D6").Value = 23999
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3

# ID Venta=6, ABC019, 2023-06-20, 49999.990000000005, 4, 3
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "ABC019"
Set-TextValue $ws.Range("C7") "2023-06-20"
$ws.Range("D7").Value = 49999.990000000005
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 3

# ID Venta=7, ABC023, 2023-07-25, 47999, 5, 2
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "ABC023"
Set-TextValue $ws.Range("C8") "2023-07-25"
$ws.Range("D8").Value = 47999
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 2

# New row 9: ID Venta=8, ABC027, 2023-08-30, 54999, 6, 1
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "ABC027"
Set-TextValue $ws.Range("C9") "2023-08-30"
$ws.Range("D9").Value = 54999
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 1

Write-Host "Applied dashboard data refresh: 8 sales rows (A2:F9)"
